# updated main GSC export data
#
# Appends the next three days of GSC export data (2025-11-22, 2025-11-23,
# 2025-11-24) as new rows 48-50 below the existing last row (47) on the
# "Chart" sheet: Date | Non-HTTPS URLs | HTTPS URLs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# The date column is plain text (e.g. "2025-11-21") in the existing rows,
# not an Excel date serial. Forcing the number format to Text before the
# assignment (and clearing it back afterwards) keeps these new cells as
# literal strings instead of being auto-parsed into dates.
$ws.Range("A48:A50").NumberFormat = "@"
$ws.Cells.Item(48, 1).Value = "2025-11-22"
$ws.Cells.Item(49, 1).Value = "2025-11-23"
$ws.Cells.Item(50, 1).Value = "2025-11-24"
$ws.Range("A48:A50").ClearFormats()

# Non-HTTPS URLs column - still 0 for every day.
$ws.Cells.Item(48, 2).Value = 0
$ws.Cells.Item(49, 2).Value = 0
$ws.Cells.Item(50, 2).Value = 0

# HTTPS URLs column.
$ws.Cells.Item(48, 3).Value = 26
$ws.Cells.Item(49, 3).Value = 26
$ws.Cells.Item(50, 3).Value = 25
